$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update event name in B8: "So Fresh Concert Autumn 2004" -> "So Fresh Concert"
$ws.Range("B8").Value = "So Fresh Concert"

# Update image filename in I11: "wildcatssydneykings.jpg" -> "wildcatskings.jpg"
$ws.Range("I11").Value = "wildcatskings.jpg"

# Update selected cell on the active sheet view
$ws.Range("I12").Select()
